$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A152").Value = 1582761600
$ws.Range("B152").NumberFormat = "@"
$ws.Range("B152").Value = "2020-02-27"
$ws.Range("C152").NumberFormat = "@"
$ws.Range("C152").Value = "0210"
$ws.Range("D152").Value = "KHJB"
$ws.Range("E152").Value = 0.255
$ws.Range("F152").Value = 0.255
$ws.Range("G152").Value = 0.255
$ws.Range("H152").Value = 0.255
$ws.Range("I152").Value = "-"

$ws.Range("A153").Value = 1582848000
$ws.Range("B153").NumberFormat = "@"
$ws.Range("B153").Value = "2020-02-28"
$ws.Range("C153").NumberFormat = "@"
$ws.Range("C153").Value = "0210"
$ws.Range("D153").Value = "KHJB"
$ws.Range("E153").Value = 0.25
$ws.Range("F153").Value = 0.26
$ws.Range("G153").Value = 0.24
$ws.Range("H153").Value = 0.245
$ws.Range("I153").Value = 415500
